# Program Flow Control.pptx - "update Programming Basics tutorial"
#
# The deck had a duplicated "switch statement" intro slide: slide 14
# ("The break statement", full content) sat right before slide 15
# ("The switch statement", a short duplicate intro). The fix removes
# the stray slide 14 so the deck flows break-content -> (deleted) ->
# switch intro -> syntax -> flowcharts -> example, i.e. slide 14 is
# simply deleted and everything else shifts up by one.
#
# It also bumps the cached "date last edited" field on the Notes
# Master from 9/18/22 to 9/19/22.

$p = $ppt.ActivePresentation

# Remove the stray "The break statement" slide (was slide 14).
$p.Slides.Item(14).Delete()

# Update the fixed date shown on notes pages (Notes Master date
# placeholder) from 9/18/22 to 9/19/22.
$nm = $p.NotesMaster
$dateAndTime = $nm.HeadersFooters.DateAndTime
$dateAndTime.Text = "9/19/22"
